# "Address check now using forms"
# Replace the old combined SeparateProducts alternate-regex cell (row 54, col C)
# with two new named regex entries (rows 54-55, cols A/B) used by the new
# address-label form-based matching: AddressContinueRegex / AddressRegex.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Order of assignment matters: it controls the insertion order into the
# shared-string table, which must come out as:
#   193 = AddressContinueRegex
#   194 = AddressRegex
#   195 = \WContinue\W:(\w+)
#   196 = {\WIN_strAddressLabel\W:\W([^(},{)]+)
$ws.Cells.Item(54, 1).Value = "AddressContinueRegex"
$ws.Cells.Item(55, 1).Value = "AddressRegex"
$ws.Cells.Item(54, 2).Value = "\WContinue\W:(\w+)"
$ws.Cells.Item(55, 2).Value = "{\WIN_strAddressLabel\W:\W([^(},{)]+)"

# Row 54 previously held its lone value in column C (with a 45pt row height
# sized for the old long wrapped regex string). Remove that cell entirely
# and let the row height fall back to the sheet default.
$ws.Cells.Item(54, 3).Clear()
$ws.Rows.Item(54).EntireRow.AutoFit()

# New blank row 55 was introduced below, pushing the used range down by one
# row and moving the active selection.
$ws.Range("B64").Select()
